$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.431.84"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "2.597.40"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "591.79"
$ws.Range("E5").Value = "  -2.47%  "
$ws.Range("D6").Value = "150.80"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D11").Value = "0.383"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "27.48"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "3.065.85"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").Value = "63.308.99"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("E16").Value = "  +6.97%  "
$ws.Range("D17").Value = "2.609.12"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "12.34"
$ws.Range("E18").Value = "  +6.30%  "
$ws.Range("D19").Value = "4.72"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("D20").Value = "345.39"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").Value = "6.91"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "67.68"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("D24").Value = "1.69"
$ws.Range("E24").Value = "  +4.69%  "
$ws.Range("D25").Value = "9.27"
$ws.Range("E25").Value = "  +2.50%  "
$ws.Range("D26").Value = "1.68"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "557.47"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").Value = "8.02"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").Value = "0.160"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").Value = "2.05"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").Value = "0.0₃0847"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("D34").Value = "5.18"
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("D35").Value = "167.23"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("D36").Value = "0.414"
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "19.58"
$ws.Range("E38").Value = "  +3.17%  "
$ws.Range("D39").Value = "1.93"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "168.05"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").Value = "3.95"
$ws.Range("E43").Value = "  +4.71%  "
$ws.Range("D44").Value = "0.0585"
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("D45").Value = "22.17"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("D47").Value = "0.0253"
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("D48").Value = "2.04"
$ws.Range("E48").Value = "  +3.75%  "
$ws.Range("D49").Value = "0.0963"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").Value = "19.12"
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("D51").Value = "0.0₆0232"
$ws.Range("E51").Value = "  +18.19%  "
